$wb = $excel.ActiveWorkbook
$wsReco = $wb.Worksheets.Item("Recommandations")
$wsYtd = $wb.Worksheets.Item("Top_YTD")

# --- Sheet "Recommandations": update changed cells (rows 2-44) ---
$wsReco.Range("D2").Value = 3635
$wsReco.Range("E2").Value = 905
$wsReco.Range("D3").Value = 3421.09
$wsReco.Range("A4").Value = 'AIR LIQUIDE CI'
$wsReco.Range("C4").Value = 4
$wsReco.Range("D4").Value = 2785
$wsReco.Range("E4").Value = 665
$wsReco.Range("A5").Value = 'BRVM - AUTRES SECTEURS'
$wsReco.Range("D5").Value = 2348.42
$wsReco.Range("E5").Value = 590.14
$wsReco.Range("A6").Value = 'BRVM - DISTRIBUTION'
$wsReco.Range("D6").Value = 2112.57
$wsReco.Range("E6").Value = 537.36
$wsReco.Range("A7").Value = 'SUCRIVOIRE'
$wsReco.Range("C7").Value = 2
$wsReco.Range("D7").Value = 1980
$wsReco.Range("E7").Value = 990
$wsReco.Range("D8").Value = 1435.59
$wsReco.Range("E8").Value = 358.59
$wsReco.Range("D9").Value = 1369.31
$wsReco.Range("E9").Value = 342.26
$wsReco.Range("D10").Value = 730.75
$wsReco.Range("E10").Value = 184.8
$wsReco.Range("A11").Value = 'BRVM - CONSOMMATION DE BASE         (**)'
$wsReco.Range("C11").Value = 3
$wsReco.Range("D11").Value = 674.79
$wsReco.Range("E11").Value = 225.7
$wsReco.Range("A12").Value = 'BRVM-PRINCIPAL                    (**)'
$wsReco.Range("C12").Value = 3
$wsReco.Range("D12").Value = 664.01
$wsReco.Range("E12").Value = 221.95
$wsReco.Range("A13").Value = 'BRVM - FINANCES'
$wsReco.Range("D13").Value = 583.23
$wsReco.Range("E13").Value = 148.25
$wsReco.Range("A14").Value = 'BRVM - SERVICES FINANCIERS'
$wsReco.Range("D14").Value = 573.1900000000001
$wsReco.Range("E14").Value = 145.69
$wsReco.Range("A15").Value = 'BRVM-PRESTIGE'
$wsReco.Range("C15").Value = 4
$wsReco.Range("D15").Value = 569.9
$wsReco.Range("E15").Value = 146.35
$wsReco.Range("A16").Value = 'BRVM - INDUSTRIE                 (**)'
$wsReco.Range("C16").Value = 2
$wsReco.Range("D16").Value = 536.74
$wsReco.Range("E16").Value = 269.25
$wsReco.Range("A17").Value = 'BRVM - INDUSTRIELS'
$wsReco.Range("C17").Value = 4
$wsReco.Range("D17").Value = 487.38
$wsReco.Range("E17").Value = 123.14
$wsReco.Range("A18").Value = 'BRVM - ENERGIE'
$wsReco.Range("D18").Value = 450.85
$wsReco.Range("E18").Value = 114.29
$wsReco.Range("A19").Value = 'BRVM - TELECOMMUNICATIONS'
$wsReco.Range("C19").Value = 4
$wsReco.Range("D19").Value = 387.2
$wsReco.Range("E19").Value = 98.8
$wsReco.Range("A22").Value = 'ORANGE COTE D''IVOIRE (ORAC)'
$wsReco.Range("B22").Value = 2
$wsReco.Range("C22").Value = 1
$wsReco.Range("D22").Value = 7.54
$wsReco.Range("E22").Value = 6.86
$wsReco.Range("G22").Value = '👀 À surveiller'
$wsReco.Range("A23").Value = 'ECOBANK COTE D''''IVOIRE (ECOC)'
$wsReco.Range("D23").Value = 7.5
$wsReco.Range("E23").Value = 7.5
$wsReco.Range("A26").Value = 'SOCIETE GENERALE COTE D''IVOIRE (SGBC)'
$wsReco.Range("D26").Value = 4.15
$wsReco.Range("E26").Value = 4.15
$wsReco.Range("A27").Value = 'SICABLE CI (CABC)'
$wsReco.Range("B27").Value = 2
$wsReco.Range("C27").Value = 1
$wsReco.Range("D27").Value = 3.69
$wsReco.Range("E27").Value = 4.07
$wsReco.Range("G27").Value = '👀 À surveiller'
$wsReco.Range("A28").Value = 'ORAGROUP TOGO (ORGT)'
$wsReco.Range("C28").Value = 1
$wsReco.Range("D28").Value = 3.14
$wsReco.Range("E28").Value = 5.91
$wsReco.Range("G28").Value = '👀 À surveiller'
$wsReco.Range("A32").Value = 'SAPH CI (SPHC)'
$wsReco.Range("D32").Value = 2.08
$wsReco.Range("E32").Value = -1.84
$wsReco.Range("A33").Value = 'SETAO CI (STAC)'
$wsReco.Range("B33").Value = 2
$wsReco.Range("D33").Value = 0.45
$wsReco.Range("E33").Value = 5.02
$wsReco.Range("C35").Value = 2
$wsReco.Range("A36").Value = 'VIVO ENERGY CI (SHEC)'
$wsReco.Range("B36").Value = 1
$wsReco.Range("D36").Value = -0.3
$wsReco.Range("E36").Value = -1.62
$wsReco.Range("A38").Value = 'SMB CI (SMBC)'
$wsReco.Range("D38").Value = -3.11
$wsReco.Range("E38").Value = -3.11
$wsReco.Range("A39").Value = 'BICI CI (BICC)'
$wsReco.Range("D39").Value = -3.23
$wsReco.Range("E39").Value = -3.23
$wsReco.Range("A40").Value = 'ECOBANK TRANS. INCORP. TG (ETIT)'
$wsReco.Range("D40").Value = -4.35
$wsReco.Range("E40").Value = -4.35
$wsReco.Range("C41").Value = 2
$wsReco.Range("D41").Value = -5.92
$wsReco.Range("E41").Value = -2.68
$wsReco.Range("A42").Value = 'SICOR CI (SICC)'
$wsReco.Range("B42").Value = 1
$wsReco.Range("D42").Value = -6.37
$wsReco.Range("E42").Value = -6.87
$wsReco.Range("A43").Value = 'NEI-CEDA CI (NEIC)'
$wsReco.Range("B43").Value = 0
$wsReco.Range("D43").Value = -6.42
$wsReco.Range("E43").Value = -2.16
$wsReco.Range("G43").Value = '➖ Neutre'
$wsReco.Range("A44").Value = 'AIR LIQUIDE CI (SIVC)'
$wsReco.Range("D44").Value = -6.99
$wsReco.Range("E44").Value = -6.99

# --- Remove now-deleted rows 45-47 on "Recommandations" ---
$wsReco.Range("A45:G47").EntireRow.Delete()

# --- Sheet "Top_YTD": update changed cells (rows 2-11) ---
$wsYtd.Range("B2").Value = 10207307.6
$wsYtd.Range("B3").Value = 1035050
$wsYtd.Range("B4").Value = 401372
$wsYtd.Range("B5").Value = 222781.45
$wsYtd.Range("B6").Value = 155513.46
$wsYtd.Range("A7").Value = 'BRVM - TRANSPORT'
$wsYtd.Range("B7").Value = 44243.89
$wsYtd.Range("A8").Value = 'BRVM - AGRICULTURE'
$wsYtd.Range("B8").Value = 38178.83
$wsYtd.Range("A9").Value = 'SUCRIVOIRE'
$wsYtd.Range("B9").Value = 11781
$wsYtd.Range("B10").Value = 6284.14
$wsYtd.Range("B11").Value = 3550.49

Write-Host "BRVM update applied"
